# Adds a "ChangePassword" module (worksheet) to the Login test-data workbook,
# matching the "Added change password module" commit.

$wb = $excel.ActiveWorkbook

# --- Create the new sheet as the LAST tab (after AddAlbum) -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ChangePassword"

# --- Headers + sample data ---------------------------------------------
# Written in this order so the shared-string table gets populated in the
# same sequence as the source workbook (New Password, Old Password,
# Confirm Password, Test@1234, Test@12345).
$ws.Range("B1").Value = "New Password"
$ws.Range("A1").Value = "Old Password"
$ws.Range("C1").Value = "Confirm Password"
$ws.Range("B2").Value = "Test@1234"
$ws.Range("A2").Value = "Test@12345"
$ws.Range("C2").Value = "Test@1234"

# --- Hyperlinks on the sample values (mirrors Login sheet's mailto style) --
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Test@1234")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:Test@12345")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Test@1234")

# --- Fit the columns to their content ---------------------------------
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# --- Page setup (portrait, like the other sheets) ----------------------
$ws.PageSetup.Orientation = 1

# --- Final selection on the new (now active) sheet ----------------------
$ws.Range("E2").Select() | Out-Null
